$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume snapshot cells (GitHub Actions refresh).
# Price cells in column D store text (e.g. "37.805.14" uses dots as both
# thousands separators and decimal separator, so Excel can't treat them
# as numbers). For values that *do* look like plain numbers (e.g. "58.76")
# a leading apostrophe forces Excel to keep them as text, matching the
# original cell type, instead of silently converting them to numeric.

$ws.Range("D2").Value = "'37.800.28"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "'2.085.70"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'234.13"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'58.76"
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").Value = "'2.394.16"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "'14.78"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "'21.22"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").Value = "'0.769"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").Value = "'2.085.40"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'37.702.35"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'6.25"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'71.48"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "'0.0₃0830"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'228.87"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "'170.48"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  +8.10%  "
$ws.Range("D28").Value = "'9.05"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("E30").Value = "  +2.49%  "
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "'3.47"
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "'1.83"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -3.84%  "
$ws.Range("E40").Value = "  +2.18%  "

# Rows 41/42 swapped position (HuobiToken now ranks above Aave) with fresh
# price/volume figures.
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'98.78"
$ws.Range("E42").Value = "  +1.07%  "

$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "'1.458.88"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").Value = "'16.12"
$ws.Range("E47").Value = "  +6.70%  "
$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").Value = "'7.45"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "'47.33"
$ws.Range("E51").Value = "  +4.90%  "
